# Eldata.xlsx — finalize production module KL-E
# Adds two new worksheets ("elasKL-E" and "prodKL-E") that extend the
# existing "elasKL" / "prodKL" sheets with an extra ENER(gy) related
# column, and makes "prodKL-E" the active sheet/tab.

$wb = $excel.ActiveWorkbook

$wsElasKL = $wb.Worksheets.Item("elasKL")
$wsProdKL = $wb.Worksheets.Item("prodKL")

# --- Create "elasKL-E" as a copy of "elasKL", placed after "prodKL" ---
$wsElasKL.Copy([System.Reflection.Missing]::Value, $wsProdKL)
$wsElasKLE = $wb.Worksheets.Item($wsProdKL.Index + 1)
$wsElasKLE.Name = "elasKL-E"

# elasKL only has an "industry"/"elasKL" header (cols A:B); elasKL-E adds
# an extra "elasKLE" column before it and an "elasE" column after it.
$wsElasKLE.Range("B1").Value = "elasKLE"
$wsElasKLE.Range("C1").Value = "elasKL"
$wsElasKLE.Range("D1").Value = "elasE"

# --- Create "prodKL-E" as a copy of "prodKL", placed after "elasKL-E" ---
$wsProdKL.Copy([System.Reflection.Missing]::Value, $wsElasKLE)
$wsProdKLE = $wb.Worksheets.Item($wsElasKLE.Index + 1)
$wsProdKLE.Name = "prodKL-E"

# prodKL has "industry"/"COE"/"GOS" headers (cols A:C); prodKL-E adds an
# extra "ENER" column.
$wsProdKLE.Range("D1").Value = "ENER"

# "prodKL-E" becomes the active sheet/tab.
$wsProdKLE.Activate()
